$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 74 values
$ws.Range("B74").Value = 0.8
$ws.Range("D74").Value = 74.5
$ws.Range("E74").Value = -9.5

# Add new row 75
$ws.Range("A75").NumberFormat = "@"
$ws.Range("A75").Value = "01-04-2021"
$ws.Range("B75").Value = -1.1
$ws.Range("C75").Value = 15.2
$ws.Range("D75").Value = 72.1
$ws.Range("E75").Value = -5.6
